$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 29 changes
$ws.Range("J2").Value = 1.95
$ws.Range("K2").Value = 2.38
$ws.Range("L2").Value = 7
$ws.Range("N2").Value = 12
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.86
$ws.Range("R2").Value = 2.04
$ws.Range("S2").Value = 1.36
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 2.05
$ws.Range("V2").Value = 1.7
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 9.5
$ws.Range("AB2").Value = 29
$ws.Range("AC2").Value = 12
$ws.Range("AD2").Value = 9
$ws.Range("AG2").Value = 451
$ws.Range("AH2").Value = 17
$ws.Range("AI2").Value = 34
$ws.Range("AJ2").Value = 21
$ws.Range("AN2").Value = 3.4
$ws.Range("AO2").Value = 7
$ws.Range("AQ2").Value = 21
$ws.Range("AS2").Value = 151
$ws.Range("AT2").Value = 3
$ws.Range("AV2").Value = 67
$ws.Range("AW2").Value = 8
$ws.Range("BC2").Value = 151

# Row 3: 1 changes
$ws.Range("N3").Value = 17

# Row 4: 26 changes
$ws.Range("G4").Value = 1.42
$ws.Range("I4").Value = 6.5
$ws.Range("J4").Value = 1.91
$ws.Range("L4").Value = 6
$ws.Range("S4").Value = 1.29
$ws.Range("T4").Value = 3.5
$ws.Range("U4").Value = 1.73
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = 9
$ws.Range("X4").Value = 8
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 10
$ws.Range("AH4").Value = 19
$ws.Range("AI4").Value = 34
$ws.Range("AJ4").Value = 19
$ws.Range("AK4").Value = 67
$ws.Range("AL4").Value = 41
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 7
$ws.Range("AQ4").Value = 19
$ws.Range("AT4").Value = 3.5
$ws.Range("AU4").Value = 8
$ws.Range("AW4").Value = 8
$ws.Range("AX4").Value = 29
$ws.Range("AZ4").Value = 101
$ws.Range("BA4").Value = 101

# Row 9: 29 changes
$ws.Range("G9").Value = 3.3
$ws.Range("I9").Value = 2.38
$ws.Range("J9").Value = 3.75
$ws.Range("L9").Value = 3.1
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
$ws.Range("Q9").Value = 2.2
$ws.Range("R9").Value = 1.65
$ws.Range("W9").Value = 9
$ws.Range("X9").Value = 15
$ws.Range("Y9").Value = 12
$ws.Range("Z9").Value = 34
$ws.Range("AA9").Value = 29
$ws.Range("AD9").Value = 5.5
$ws.Range("AH9").Value = 7.5
$ws.Range("AI9").Value = 11
$ws.Range("AJ9").Value = 9.5
$ws.Range("AK9").Value = 23
$ws.Range("AL9").Value = 21
$ws.Range("AM9").Value = 34
$ws.Range("AN9").Value = 5
$ws.Range("AO9").Value = 17
$ws.Range("AQ9").Value = 51
$ws.Range("AR9").Value = 81
$ws.Range("AS9").Value = 201
$ws.Range("AU9").Value = 8
$ws.Range("AW9").Value = 4.33
$ws.Range("AX9").Value = 13
$ws.Range("AY9").Value = 26

# Row 10: 34 changes
$ws.Range("G10").Value = 4.1
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 1.91
$ws.Range("J10").Value = 4.33
$ws.Range("L10").Value = 2.6
$ws.Range("O10").Value = 1.25
$ws.Range("P10").Value = 3.75
$ws.Range("Q10").Value = 1.9
$ws.Range("R10").Value = 1.95
$ws.Range("S10").Value = 1.36
$ws.Range("T10").Value = 3
$ws.Range("U10").Value = 1.75
$ws.Range("V10").Value = 2
$ws.Range("W10").Value = 12
$ws.Range("X10").Value = 21
$ws.Range("Y10").Value = 13
$ws.Range("Z10").Value = 41
$ws.Range("AA10").Value = 34
$ws.Range("AC10").Value = 11
$ws.Range("AD10").Value = 6.5
$ws.Range("AI10").Value = 9.5
$ws.Range("AJ10").Value = 8.5
$ws.Range("AK10").Value = 17
$ws.Range("AL10").Value = 15
$ws.Range("AO10").Value = 21
$ws.Range("AP10").Value = 29
$ws.Range("AQ10").Value = 67
$ws.Range("AS10").Value = 201
$ws.Range("AT10").Value = 3
$ws.Range("AU10").Value = 8
$ws.Range("AW10").Value = 4
$ws.Range("AX10").Value = 10
$ws.Range("AZ10").Value = 34
$ws.Range("BB10").Value = 126

# Row 11: 28 changes
$ws.Range("G11").Value = 2.15
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 3.5
$ws.Range("J11").Value = 2.88
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 4.33
$ws.Range("Q11").Value = 2.35
$ws.Range("R11").Value = 1.57
$ws.Range("W11").Value = 6.5
$ws.Range("X11").Value = 9.5
$ws.Range("Y11").Value = 9.5
$ws.Range("Z11").Value = 19
$ws.Range("AA11").Value = 19
$ws.Range("AD11").Value = 6
$ws.Range("AG11").Value = 401
$ws.Range("AH11").Value = 8.5
$ws.Range("AI11").Value = 17
$ws.Range("AJ11").Value = 13
$ws.Range("AL11").Value = 34
$ws.Range("AN11").Value = 4
$ws.Range("AO11").Value = 12
$ws.Range("AP11").Value = 26
$ws.Range("AU11").Value = 8.5
$ws.Range("AW11").Value = 5.5
$ws.Range("AX11").Value = 21
$ws.Range("AZ11").Value = 67
$ws.Range("BA11").Value = 101
$ws.Range("BB11").Value = 251

# Row 12: 33 changes
$ws.Range("G12").Value = 1.2
$ws.Range("H12").Value = 7
$ws.Range("I12").Value = 13
$ws.Range("J12").Value = 1.53
$ws.Range("K12").Value = 3.1
$ws.Range("L12").Value = 9
$ws.Range("N12").Value = 23
$ws.Range("O12").Value = 1.1
$ws.Range("P12").Value = 7
$ws.Range("Q12").Value = 1.33
$ws.Range("R12").Value = 3.4
$ws.Range("U12").Value = 1.8
$ws.Range("V12").Value = 1.95
$ws.Range("Y12").Value = 10
$ws.Range("AB12").Value = 23
$ws.Range("AC12").Value = 23
$ws.Range("AD12").Value = 13
$ws.Range("AE12").Value = 21
$ws.Range("AG12").Value = 201
$ws.Range("AJ12").Value = 34
$ws.Range("AK12").Value = 151
$ws.Range("AL12").Value = 67
$ws.Range("AM12").Value = 51
$ws.Range("AN12").Value = 3.6
$ws.Range("AO12").Value = 5.5
$ws.Range("AQ12").Value = 11
$ws.Range("AR12").Value = 29
$ws.Range("AU12").Value = 9
$ws.Range("AW12").Value = 12
$ws.Range("AX12").Value = 41
$ws.Range("AZ12").Value = 201
$ws.Range("BB12").Value = 201
$ws.Range("BC12").Value = 451

# Row 13: 20 changes
$ws.Range("G13").Value = 1.91
$ws.Range("H13").Value = 3.5
$ws.Range("I13").Value = 3.9
$ws.Range("J13").Value = 2.6
$ws.Range("K13").Value = 2.2
$ws.Range("L13").Value = 4.33
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 11
$ws.Range("X13").Value = 9
$ws.Range("Y13").Value = 8.5
$ws.Range("Z13").Value = 17
$ws.Range("AA13").Value = 15
$ws.Range("AC13").Value = 11
$ws.Range("AH13").Value = 11
$ws.Range("AI13").Value = 21
$ws.Range("AJ13").Value = 13
$ws.Range("AO13").Value = 10
$ws.Range("AQ13").Value = 34
$ws.Range("AX13").Value = 21
$ws.Range("AY13").Value = 29

# Row 14: 2 changes
$ws.Range("Q14").Value = 1.93
$ws.Range("R14").Value = 1.93

# Row 15: 10 changes
$ws.Range("G15").Value = 3.4
$ws.Range("H15").Value = 3.25
$ws.Range("I15").Value = 2.15
$ws.Range("K15").Value = 2.1
$ws.Range("N15").Value = 8.5
$ws.Range("Z15").Value = 41
$ws.Range("AJ15").Value = 9
$ws.Range("AK15").Value = 19
$ws.Range("AN15").Value = 5.5
$ws.Range("BB15").Value = 151

# Row 19: 11 changes
$ws.Range("G19").Value = 1.42
$ws.Range("K19").Value = 2.75
$ws.Range("Q19").Value = 1.4
$ws.Range("R19").Value = 2.88
$ws.Range("U19").Value = 1.62
$ws.Range("V19").Value = 2.2
$ws.Range("W19").Value = 10
$ws.Range("Y19").Value = 8.5
$ws.Range("AK19").Value = 67
$ws.Range("AU19").Value = 8
$ws.Range("AY19").Value = 29

# Row 20: 24 changes
$ws.Range("G20").Value = 2.2
$ws.Range("H20").Value = 3.6
$ws.Range("I20").Value = 3.1
$ws.Range("J20").Value = 2.75
$ws.Range("K20").Value = 2.38
$ws.Range("L20").Value = 3.4
$ws.Range("O20").Value = 1.17
$ws.Range("P20").Value = 5
$ws.Range("Q20").Value = 1.6
$ws.Range("R20").Value = 2.3
$ws.Range("Y20").Value = 9
$ws.Range("Z20").Value = 21
$ws.Range("AA20").Value = 15
$ws.Range("AE20").Value = 12
$ws.Range("AG20").Value = 101
$ws.Range("AH20").Value = 13
$ws.Range("AI20").Value = 19
$ws.Range("AK20").Value = 34
$ws.Range("AM20").Value = 26
$ws.Range("AN20").Value = 4.5
$ws.Range("AO20").Value = 11
$ws.Range("AQ20").Value = 34
$ws.Range("AW20").Value = 5.5
$ws.Range("AZ20").Value = 51

Write-Host "Applied 247 cell updates"